$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.603.64"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.640.33"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'536.97"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'146.13"
$ws.Range("E6").Value = "  +3.63%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "'6.89"
$ws.Range("E9").Value = "  +6.88%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "3.110.74"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "59.500.29"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("D15").Value = "'21.39"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "2.666.16"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "'4.49"
$ws.Range("E18").Value = "  +2.98%  "
$ws.Range("D19").Value = "'339.64"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'10.32"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").Value = "'6.21"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'66.32"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").Value = "'0.417"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'7.31"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("D28").Value = "0.0₃0748"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").Value = "'5.86"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").Value = "'18.86"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").Value = "'151.17"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("D37").Value = "'0.838"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").Value = "'3.61"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").Value = "'285.38"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'10.74"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0540"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("D46").Value = "'0.0945"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "1.962.29"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.57"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'18.45"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'111.48"
$ws.Range("E51").Value = "  +0.38%  "
